$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("searchDetails")
$ws.Activate()

# New header for the additional expected-result column
$ws.Range("C1").Value = "searchValue2"

# Existing rows keep their search-type / search-value but gain an expected
# "searchValue2" (second expected value) column
$ws.Range("C2").Value = "null"
$ws.Range("C3").Value = "null"
$ws.Range("C4").Value = "null"

# The two separate zone-number / ward-number search rows are merged into a
# single "search with zone and ward number" scenario with combined values
$ws.Range("A5").Value = "searchWithZoneAndWardNumber"
$ws.Range("B5").Value = "Zone-15;Revenue Ward No  87"
$ws.Range("C5").Value = "87/1110-9-c;C. Naga Sailaja W/o R. Satish Kumar"

# New scenario: search by owner name
$ws.Range("A6").Value = "searchWithOwnerName"
$ws.Range("B6").Value = "revenue colony;C. Naga Sailaja W/o R. Satish Kumar"

# New scenario: search by demand
$ws.Range("A7").Value = "searchByDemand"
$ws.Range("B7").Value = "500;501"

$ws.Columns.Item(2).ColumnWidth = 43.33
$ws.Columns.Item(3).ColumnWidth = 40.09

$ws.Range("B7").Select()
